$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price column D, Volume(1h) column E)
$data = @{
  2 = @("30.652.47", "  +0.48%  ")
  3 = @("2.115.70", "  +0.28%  ")
  4 = @("1.014", "  +1.12%  ")
  5 = @("338.85", "  +1.38%  ")
  6 = @("1.012", "  +1.00%  ")
  7 = @("0.5254", "  +0.21%  ")
  8 = @("0.4535", "  +0.20%  ")
  9 = @("54.16", "  +1.31%  ")
  10 = @("0.09089", "  +0.99%  ")
  11 = @("1.172", "  +0.52%  ")
  12 = @("24.39", "  -0.42%  ")
  13 = @("2.122.68", "  +0.93%  ")
  14 = @("6.819", "  +0.51%  ")
  15 = @("8.084", "  +3.39%  ")
  16 = @("97.71", "  +1.01%  ")
  17 = @($null, "  +3.26%  ")
  18 = @("1.014", "  +1.04%  ")
  19 = @("0.06705", "  +1.19%  ")
  20 = @("19.36", "  +0.23%  ")
  21 = @($null, "  +1.07%  ")
  22 = @("6.410", "  +1.63%  ")
  23 = @("30.759.07", "  +0.64%  ")
  24 = @($null, "  +3.63%  ")
  25 = @("2.373", "  +1.32%  ")
  26 = @("2.369.88", "  +0.85%  ")
  27 = @("22.40", $null)
  28 = @("165.04", "  +0.90%  ")
  29 = @("2.549", "  -1.29%  ")
  30 = @("135.67", "  +2.16%  ")
  31 = @("1.199", "  -0.30%  ")
  32 = @($null, "  +0.19%  ")
  33 = @("6.375", "  +3.27%  ")
  34 = @("1.635", "  -1.56%  ")
  35 = @("3.945", "  +0.06%  ")
  36 = @("10.36", "  -3.07%  ")
  37 = @("5.919", "  +6.70%  ")
  38 = @("0.02653", "  +2.83%  ")
  39 = @("0.06838", "  +0.06%  ")
  40 = @($null, "  +1.39%  ")
  41 = @($null, "  -1.53%  ")
  42 = @("0.6882", "  -0.82%  ")
  43 = @("1.261", "  +1.30%  ")
  44 = @("15.09", "  +6.89%  ")
  45 = @($null, "  +0.22%  ")
  46 = @("2.317", "  -3.61%  ")
  47 = @("0.00000000370", "  +15.46%  ")
  48 = @("3.704", "  +1.19%  ")
  49 = @("1.255", "  +0.49%  ")
  50 = @("0.07320", "  +3.26%  ")
  51 = @("82.82", "  -0.53%  ")
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  $dVal = $vals[0]
  $eVal = $vals[1]
  if ($dVal -ne $null) {
    $dCell = $ws.Range("D" + $row)
    $dCell.NumberFormat = "@"
    $dCell.Value = $dVal
  }
  if ($eVal -ne $null) {
    $ws.Range("E" + $row).Value = $eVal
  }
}
